$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "Zinc inhibits intestinal copper uptake, an effect utilized for treating Wilson's disease (WD). We used copper-64 (64Cu) PET/CT to examine how much four weeks of treatment with different zinc regimens reduced the hepatic 64Cu content after oral 64Cu administration and test if alternative regimens were noninferior to the standard regimen of zinc acetate 50 mg × 3 daily. Forty healthy persons were randomized to four different zinc protocols. The WD standard treatment zinc acetate 50 mg × 3 reduced the hepatic 64Cu content from 26.9 ± 7.5% to 13.3 ± 5.6% of the administered 64Cu. Zinc gluconate 50 mg × 3 was noninferior (P = 0.02) (35.8 ± 9.0% to 17.4 ± 7.5%). Zinc acetate 150 mg × 1 (33.1 ± 9.9% to 17.4 ± 7.5%) and zinc gluconate 150 mg × 1 (28.1 ± 6.7% to 22.0 ± 6.7%) were less effective. These effects were intra- and inter-individually highly variable, and 14% had no effect of any zinc regimen, which may explain disparities in zinc treatment efficacy in WD patients."
$ws.Range("AF9").Value = "c(""https://openalex.org/W2469850045"", ""https://openalex.org/W2417131764"", ""https://openalex.org/W2070652985"", ""https://openalex.org/W2007094274"", ""https://openalex.org/W2085980434"", ""https://openalex.org/W1632137160"", ""https://openalex.org/W2072035168"", ""https://openalex.org/W1509942693"", ""https://openalex.org/W171658490"", ""https://openalex.org/W2358136759"")"
$ws.Range("C11").Value = "list(au_id = c(""https://openalex.org/A5053434417"", ""https://openalex.org/A5090100242"", ""https://openalex.org/A5005707441"", ""https://openalex.org/A5006147905"", ""https://openalex.org/A5027999373"", ""https://openalex.org/A5075763247"", ""https://openalex.org/A5015984752"", ""https://openalex.org/A5059158988"", ""https://openalex.org/A5084607692"", ""https://openalex.org/A5057470830"", ""https://openalex.org/A5036172185"", ""https://openalex.org/A5014882919"", ""https://openalex.org/A5072641596"", ""https://openalex.org/A5005923405"", `n""https://openalex.org/A5009531475"", ""https://openalex.org/A5063597334"", ""https://openalex.org/A5018746315"", ""https://openalex.org/A5055020443"", ""https://openalex.org/A5085086661""), au_display_name = c(""Xu Chu"", ""Guang Yang"", ""Alexandros Terzis"", ""Visakh Vaikuntanathan"", ""Wenkang Wang"", ""Zhouhang Li"", ""Grazia Lamanna"", ""Stephanie Fest-Santini"", ""Maurizio Santini"", ""Gianpietro Cossali"", ""Phillip M. Ligrani"", ""Bassam A. Younis"", ""Michael Crawford"", ""Peter Ott"", ""Jürgen Köhler"", ""Christian Rohde"", `n""Claus-Dieter Münz"", ""Rainer Helmig"", ""Tianshou Zhao""), au_orcid = c(""https://orcid.org/0000-0001-8456-1865"", ""https://orcid.org/0000-0002-8430-6917"", NA, ""https://orcid.org/0000-0002-1781-2334"", ""https://orcid.org/0000-0003-4134-8364"", ""https://orcid.org/0000-0003-1513-6499"", ""https://orcid.org/0000-0002-1725-9663"", ""https://orcid.org/0000-0002-0613-390X"", ""https://orcid.org/0000-0003-4652-3264"", ""https://orcid.org/0000-0003-2550-6780"", ""https://orcid.org/0000-0002-0313-9580"", ""https://orcid.org/0000-0002-3437-1219"", `nNA, ""https://orcid.org/0000-0002-3088-1983"", NA, ""https://orcid.org/0000-0001-6778-5971"", NA, ""https://orcid.org/0000-0003-2601-5377"", ""https://orcid.org/0000-0002-9598-932X""), author_position = c(""first"", ""middle"", ""middle"", ""middle"", ""middle"", ""middle"", ""middle"", ""middle"", ""middle"", ""middle"", ""middle"", ""middle"", ""middle"", ""middle"", ""middle"", ""middle"", ""middle"", ""middle"", ""last""), au_affiliation_raw = c(""University of Stuttgart, Germany"", ""Shanghai Jiao Tong University, China"", ""Technion - Israel Institute of Technology, Israel"", `n""Shiv Nadar University, India"", ""University of Stuttgart, Germany"", ""Kunming University of Science and Technology, China"", ""University of Stuttgart, Germany"", ""Università degli Studi di Bergamo, Italy"", ""Università degli Studi di Bergamo, Italy"", ""Università degli Studi di Bergamo, Italy"", ""University of Alabama at Huntsville, USA"", ""University of California - Davis, USA"", ""The University of Texas at Austin, USA"", ""EPFL, Switzerland"", ""TU Braunschweig, Germany"", ""University of Stuttgart, Germany"", `n""University of Stuttgart, Germany"", ""University of Stuttgart, Germany"", ""Hong Kong University of Science and Technology, China""), institution_id = c(""https://openalex.org/I100066346"", ""https://openalex.org/I183067930"", ""https://openalex.org/I174306211"", ""https://openalex.org/I26604189"", ""https://openalex.org/I100066346"", ""https://openalex.org/I10660446"", ""https://openalex.org/I100066346"", ""https://openalex.org/I11039511"", ""https://openalex.org/I11039511"", ""https://openalex.org/I11039511"", ""https://openalex.org/I82495205"", `n""https://openalex.org/I84218800"", ""https://openalex.org/I86519309"", ""https://openalex.org/I5124864"", ""https://openalex.org/I94509681"", ""https://openalex.org/I100066346"", ""https://openalex.org/I100066346"", ""https://openalex.org/I100066346"", ""https://openalex.org/I889458895""), institution_display_name = c(""University of Stuttgart"", ""Shanghai Jiao Tong University"", ""Technion – Israel Institute of Technology"", ""Shiv Nadar University"", ""University of Stuttgart"", ""Kunming University of Science and Technology"", `n""University of Stuttgart"", ""University of Bergamo"", ""University of Bergamo"", ""University of Bergamo"", ""University of Alabama in Huntsville"", ""University of California, Davis"", ""The University of Texas at Austin"", ""École Polytechnique Fédérale de Lausanne"", ""Technische Universität Braunschweig"", ""University of Stuttgart"", ""University of Stuttgart"", ""University of Stuttgart"", ""University of Hong Kong""), institution_ror = c(""https://ror.org/04vnq7t77"", ""https://ror.org/0220qvk04"", ""https://ror.org/03qryx823"", `n""https://ror.org/05aqahr97"", ""https://ror.org/04vnq7t77"", ""https://ror.org/00xyeez13"", ""https://ror.org/04vnq7t77"", ""https://ror.org/02mbd5571"", ""https://ror.org/02mbd5571"", ""https://ror.org/02mbd5571"", ""https://ror.org/02zsxwr40"", ""https://ror.org/05rrcem69"", ""https://ror.org/00hj54h04"", ""https://ror.org/02s376052"", ""https://ror.org/010nsgg66"", ""https://ror.org/04vnq7t77"", ""https://ror.org/04vnq7t77"", ""https://ror.org/04vnq7t77"", ""https://ror.org/02zhqgq86""), institution_country_code = c(""DE"", ""CN"", `n""IL"", ""IN"", ""DE"", ""CN"", ""DE"", ""IT"", ""IT"", ""IT"", ""US"", ""US"", ""US"", ""CH"", ""DE"", ""DE"", ""DE"", ""DE"", ""HK""), institution_type = c(""education"", ""education"", ""education"", ""education"", ""education"", ""education"", ""education"", ""education"", ""education"", ""education"", ""education"", ""education"", ""education"", ""education"", ""education"", ""education"", ""education"", ""education"", ""education""), institution_lineage = c(""https://openalex.org/I100066346"", ""https://openalex.org/I183067930"", ""https://openalex.org/I174306211"", `n""https://openalex.org/I26604189"", ""https://openalex.org/I100066346"", ""https://openalex.org/I10660446"", ""https://openalex.org/I100066346"", ""https://openalex.org/I11039511"", ""https://openalex.org/I11039511"", ""https://openalex.org/I11039511"", ""https://openalex.org/I2800507078, https://openalex.org/I82495205"", ""https://openalex.org/I2803209242, https://openalex.org/I84218800"", ""https://openalex.org/I16452829, https://openalex.org/I86519309"", ""https://openalex.org/I2799323385, https://openalex.org/I5124864"", `n""https://openalex.org/I94509681"", ""https://openalex.org/I100066346"", ""https://openalex.org/I100066346"", ""https://openalex.org/I100066346"", ""https://openalex.org/I889458895""))"
$ws.Range("C15").Value = "list(au_id = c(""https://openalex.org/A5054938817"", ""https://openalex.org/A5069717460"", ""https://openalex.org/A5005923405"", ""https://openalex.org/A5002371143"", ""https://openalex.org/A5083246247"", ""https://openalex.org/A5002664988"", ""https://openalex.org/A5036981869"", ""https://openalex.org/A5070215060""), au_display_name = c(""C.Oma Kamlin"", ""Michael L. Schilsky"", ""Peter Ott"", ""Karl Heinz Weiss"", ""Massimo Zuin"", ""Aurélia Poujois"", ""Aftab Ala"", ""Koenraad D’Hollander""), au_orcid = c(NA, ""https://orcid.org/0000-0001-9043-0554"", `n""https://orcid.org/0000-0002-3088-1983"", ""https://orcid.org/0000-0002-6336-9935"", NA, ""https://orcid.org/0000-0002-5640-7534"", NA, NA), author_position = c(""first"", ""middle"", ""middle"", ""middle"", ""middle"", ""middle"", ""middle"", ""last""), au_affiliation_raw = c(""Orphalan SA France, Medical Affairs, London, United Kingdom;"", ""Yale University, Division of Internal Medicine and Liver Unit, United States;"", ""Aarhus University Hospital, Aarhus, Denmark;"", ""Salem Medical Center, Heidelberg, Germany;"", ""University of Milan, Milan, Italy;"", `n""Fondation Hopital Rothschild, Department of Neurology, Paris;"", ""Kings College University Hospital, London, United Kingdom;"", ""International Drug Development Institute (IDDI), Belgium""), institution_id = c(NA, ""https://openalex.org/I32971472"", ""https://openalex.org/I2802335433"", ""https://openalex.org/I4210146747"", ""https://openalex.org/I189158943"", ""https://openalex.org/I4210127304"", ""https://openalex.org/I4210153400"", ""https://openalex.org/I4210101853""), institution_display_name = c(NA, ""Yale University"", `n""Aarhus University Hospital"", ""Krankenhaus Salem"", ""University of Milan"", ""Fondation de Rothschild"", ""King's College Hospital"", ""International Drug Development Institute (Belgium)""), institution_ror = c(NA, ""https://ror.org/03v76x132"", ""https://ror.org/040r8fr65"", ""https://ror.org/04c9pv832"", ""https://ror.org/00wjc7c48"", ""https://ror.org/02yfw7119"", ""https://ror.org/044nptt90"", ""https://ror.org/016dg3e07""), institution_country_code = c(NA, ""US"", ""DK"", ""DE"", ""IT"", ""FR"", ""GB"", ""BE""), institution_type = c(NA, `n""education"", ""healthcare"", ""healthcare"", ""education"", ""nonprofit"", ""healthcare"", ""company""), institution_lineage = c(NA, ""https://openalex.org/I32971472"", ""https://openalex.org/I2802335433"", ""https://openalex.org/I4210146747"", ""https://openalex.org/I189158943"", ""https://openalex.org/I4210127304"", ""https://openalex.org/I4210111135, https://openalex.org/I4210153400"", ""https://openalex.org/I4210101853""))"
$ws.Range("D15").Value = "Method: We retrospectively analysis a monocentric cohort of WD patients with liver damage from a referral center for rare diseases.Patients of the cohort were regularly assessed with a maximal interval of 6 months.Patients with a history of cirrhosis had a biannual ultrasound examination.Occurrence, treatment and outcome of primary liver tumors were analyzed.Results: We analyzed 149 WD patients of the cohort, with a mean follow up of 15 years.The mean age at diagnosis of WD was 20.1 years.Sixteen patients were diagnosed with a primary liver tumor during follow up.The mean age at onset of a liver tumor was 52.5 years, with a mean interval between diagnosis of WD and occurrence of tumor of 19.2 years.We diagnosed a hepatocellular carcinoma (HCC) for 10 patients and a cholangiocarcinoma (CCK) for 6 patients.Histological proof could be obtained for 7/10 patients with HCC and for all patients (6/6) with CCK.For the 3 other patient, diagnosis of HCC was based on the morphological aspect (hyper vascular tumor with washout) and elevated AFP markers.All patients who developed a liver tumor had evidence of cirrhosis at the diagnosis of WD.At the diagnostic of tumor, there was no cirrhosis on the pathological examination of the non-tumoral liver for 6 patients (2/10 patients with histological proof of HCC and 4/6 patients with CCK).For all patients with primary liver tumor, a treatment could be proposed.Resection was possible for 8 patients and liver transplantation for 5 patients.During follow up, six patients died from complications of tumor extension (3 from HCC and 3 from CCK).Conclusion: WD patients with history of cirrhosis are at risk of developing primary liver tumor and require morphological screening, even if there is a regression of fibrosis.The proportion of CCK (37.5%) among these tumors suggests an interest of the biopsy for diagnosis of liver cancer in WD patients."
$ws.Range("AF15").Value = "c(""https://openalex.org/W339681517"", ""https://openalex.org/W4322576098"", ""https://openalex.org/W2167689122"", ""https://openalex.org/W2071434358"", ""https://openalex.org/W2124699488"", ""https://openalex.org/W2053082293"", ""https://openalex.org/W1606669497"", ""https://openalex.org/W2518169464"", ""https://openalex.org/W2166447472"", ""https://openalex.org/W4288401684"")"
